# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder country list: Armenia now appears where Austria was (row 62)
# and Austria where Armenia was (row 63) - i.e. Armenia moves ahead of
# Austria in the country table.
$ws.Range("A62").Value = "Armenia"
$ws.Range("A63").Value = "Austria"

# Row 27 - Israel
$ws.Cells.Item(27, 2).Value = 299253
$ws.Cells.Item(27, 3).Value = 753
$ws.Cells.Item(27, 4).Value = 253212
$ws.Cells.Item(27, 5).Value = 43942
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 2099

# Row 28 - Ucrania
$ws.Cells.Item(28, 2).Value = 281239
$ws.Cells.Item(28, 3).Value = 5062
$ws.Cells.Item(28, 4).Value = 119650
$ws.Cells.Item(28, 5).Value = 156287
$ws.Cells.Item(28, 7).Value = 73
$ws.Cells.Item(28, 8).Value = 5302

# Row 62 - now Armenia
$ws.Cells.Item(62, 2).Value = 59995
$ws.Cells.Item(62, 3).Value = 1371
$ws.Cells.Item(62, 4).Value = 47119
$ws.Cells.Item(62, 5).Value = 11830
$ws.Cells.Item(62, 7).Value = 7
$ws.Cells.Item(62, 8).Value = 1046

# Row 63 - now Austria
$ws.Cells.Item(63, 2).Value = 58672
$ws.Cells.Item(63, 4).Value = 45846
$ws.Cells.Item(63, 5).Value = 11954
$ws.Cells.Item(63, 8).Value = 872

# Row 75 - Hungria
$ws.Cells.Item(75, 2).Value = 41732
$ws.Cells.Item(75, 3).Value = 950
$ws.Cells.Item(75, 4).Value = 12628
$ws.Cells.Item(75, 5).Value = 28052
$ws.Cells.Item(75, 7).Value = 29
$ws.Cells.Item(75, 8).Value = 1052

# Row 76 - Afganistan
$ws.Cells.Item(76, 2).Value = 40026
$ws.Cells.Item(76, 3).Value = 32
$ws.Cells.Item(76, 4).Value = 33447
$ws.Cells.Item(76, 5).Value = 5098

# Row 176 - Taiwan
$ws.Cells.Item(176, 2).Value = 531
$ws.Cells.Item(176, 3).Value = 1
$ws.Cells.Item(176, 5).Value = 33
